$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: style copied from the existing header cell (H1) so new headers match formatting
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I and J
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 10

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 6
